$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 40 updates
$ws.Range("J40").Value = 1.02
$ws.Range("K40").Value = 11
$ws.Range("L40").Value = 1.11
$ws.Range("M40").Value = 6

# Row 46 updates
$ws.Range("G46").Value = 2.42
$ws.Range("H46").Value = 3.35
$ws.Range("I46").Value = 2.7
$ws.Range("K46").Value = 7.7
$ws.Range("L46").Value = 1.27
$ws.Range("M46").Value = 3.45
$ws.Range("N46").Value = 1.8
$ws.Range("O46").Value = 1.91
$ws.Range("P46").Value = 1.38
$ws.Range("Q46").Value = 2.8
$ws.Range("R46").Value = 1.65
$ws.Range("S46").Value = 2.12
$ws.Range("T46").Value = 8.75
$ws.Range("U46").Value = 12.5
$ws.Range("W46").Value = 26
$ws.Range("X46").Value = 19
$ws.Range("Y46").Value = 27
$ws.Range("Z46").Value = 7.7
$ws.Range("AA46").Value = 6.5
$ws.Range("AD46").Value = 350
$ws.Range("AH46").Value = 32
$ws.Range("AI46").Value = 21
$ws.Range("AJ46").Value = 27
